$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows mirroring the existing table format
$ws.Range("A18").Value = "/ngo/:email"
$ws.Range("B18").Value = "delete"
$ws.Range("A19").Value = "/user/:email"
$ws.Range("B19").Value = "delete"

$ws.Rows("18:19").RowHeight = 15.75

# Match the font formatting used by column A of the other data rows (size 12 Calibri)
$ws.Range("A18").Font.Size = 12
$ws.Range("A19").Font.Size = 12

$ws.Range("A19").Select()
